$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.634463310241699
$ws.Range("B1").Value = 6.267955780029297
$ws.Range("C1").Value = 8.706923484802246
$ws.Range("D1").Value = 9.179832458496094
$ws.Range("E1").Value = 1.756473541259766
